$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("compute_instance")

# Clear the placeholder "Values" (column D) for rows 8 through 19 -
# these held generic val8..val17/val19 sample text (and the default
# network_interface JSON on row 8) that should no longer be populated.
$ws.Range("D8:D19").ClearContents()

# Add a new blank worksheet named "Sheet1" right after "compute_instance"
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "Sheet1"
